$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record for Ajo (Chino, Primera) at Feria Lagunitas de Puerto
# Montt needs to be inserted as row 388, pushing all existing records (rows
# 388:451) down by one row to (389:452). Insert a fresh row at 388 first so
# the rest of the table shifts down intact.
$ws.Rows("388:388").Insert()

# Duplicate the row that landed at 389 (the former row 388) back into the new
# row 388 so every column that doesn't change keeps its original content,
# formatting, and style.
$ws.Range("A389:R389").Copy($ws.Range("A388"))

# Now overwrite the cells that hold the new record's figures.
$ws.Range("D388").Value = 45034
$ws.Range("K388").Value = 18500
$ws.Range("L388").Value = 21000
$ws.Range("M388").Value = 19750
$ws.Range("P388").Value = 1975
